# Update "想去人数" (F) counts and one "最低票价" (G) status label
# across the four worksheets, per the scraped commit diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览 (sheet1) ----
$ws1.Range("F2").Value = 1542
$ws1.Range("F3").Value = 1505
$ws1.Range("F6").Value = 780
$ws1.Range("F7").Value = 52
$ws1.Range("F8").Value = 705
$ws1.Range("F11").Value = 1413
$ws1.Range("F12").Value = 37084
$ws1.Range("G12").Value = "已售罄"
$ws1.Range("F13").Value = 7462
$ws1.Range("F14").Value = 126
$ws1.Range("F15").Value = 410
$ws1.Range("F16").Value = 606
$ws1.Range("F18").Value = 47
$ws1.Range("F19").Value = 131
$ws1.Range("F20").Value = 471
$ws1.Range("F21").Value = 8
$ws1.Range("F23").Value = 477
$ws1.Range("F25").Value = 868
$ws1.Range("F26").Value = 39
$ws1.Range("F27").Value = 338
$ws1.Range("F28").Value = 420
$ws1.Range("F31").Value = 265
$ws1.Range("F33").Value = 761
$ws1.Range("F34").Value = 309
$ws1.Range("F35").Value = 143
$ws1.Range("F36").Value = 122
$ws1.Range("F37").Value = 795
$ws1.Range("F38").Value = 129
$ws1.Range("F40").Value = 846
$ws1.Range("F41").Value = 309

# ---- 演出 (sheet2) ----
$ws2.Range("F2").Value = 1256
$ws2.Range("F7").Value = 4336
$ws2.Range("F10").Value = 6
$ws2.Range("F12").Value = 61
$ws2.Range("F13").Value = 47
$ws2.Range("F19").Value = 4321

# ---- 本地生活 (sheet3) ----
$ws3.Range("F2").Value = 1546
$ws3.Range("F3").Value = 381

# ---- 全部类型 (sheet4) ----
$ws4.Range("F2").Value = 1546
$ws4.Range("F3").Value = 381
$ws4.Range("F4").Value = 1256
$ws4.Range("F5").Value = 1542
$ws4.Range("F7").Value = 1505
$ws4.Range("F9").Value = 780
$ws4.Range("F10").Value = 52
$ws4.Range("F11").Value = 705
$ws4.Range("F13").Value = 37084
$ws4.Range("F17").Value = 6
$ws4.Range("F19").Value = 7462
$ws4.Range("F20").Value = 410
$ws4.Range("F21").Value = 61
$ws4.Range("F22").Value = 606
$ws4.Range("F24").Value = 47
$ws4.Range("F25").Value = 47
$ws4.Range("F26").Value = 131
$ws4.Range("F27").Value = 471
$ws4.Range("F29").Value = 8
$ws4.Range("F33").Value = 868
$ws4.Range("F34").Value = 39
$ws4.Range("F35").Value = 338
$ws4.Range("F36").Value = 420
$ws4.Range("F39").Value = 265
$ws4.Range("F41").Value = 761
$ws4.Range("F43").Value = 309
$ws4.Range("F44").Value = 143
$ws4.Range("F45").Value = 122
$ws4.Range("F46").Value = 846
$ws4.Range("F47").Value = 309
